# Update F-column "views" counts across the four worksheets.
# Sheet names (workbook tab order): 展览, 演出, 本地生活, 全部类型

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 341
$ws1.Range("F4").Value = 1328
$ws1.Range("F6").Value = 363
$ws1.Range("F8").Value = 241
$ws1.Range("F9").Value = 784
$ws1.Range("F10").Value = 2356
$ws1.Range("F15").Value = 207
$ws1.Range("F16").Value = 196
$ws1.Range("F17").Value = 2747
$ws1.Range("F21").Value = 349
$ws1.Range("F22").Value = 238

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 30
$ws2.Range("F22").Value = 77

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 2125
$ws3.Range("F6").Value = 17

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2125
$ws4.Range("F10").Value = 341
$ws4.Range("F11").Value = 1328
$ws4.Range("F13").Value = 30
$ws4.Range("F16").Value = 17
$ws4.Range("F17").Value = 363
$ws4.Range("F20").Value = 241
$ws4.Range("F24").Value = 784
$ws4.Range("F25").Value = 2356
$ws4.Range("F31").Value = 207
$ws4.Range("F32").Value = 196
$ws4.Range("F40").Value = 349
$ws4.Range("F41").Value = 238
$ws4.Range("F49").Value = 77
